# Weekly update: two new price rows for Cilantro were added to the
# "Vega Central Mapocho de Santiago" sheet, pushing the old rows 896-938
# down by two (to 898-940) and growing the used range from A1:R938 to
# A1:R940.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows right before the old row 896 - this shifts all
# rows from 896 downward (896..938) down to 898..940, carrying their
# existing content (and styles) with them.
$ws.Rows.Item(896).Insert()
$ws.Rows.Item(897).Insert()

# Populate the first new row (896) with the new record.
$ws.Cells.Item(896, 1).Value = 9
$ws.Cells.Item(896, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(896, 3).Value = "Metropolitana"
$ws.Cells.Item(896, 4).Value = 45041
$ws.Cells.Item(896, 5).Value = 13
$ws.Cells.Item(896, 6).Value = 100112040
$ws.Cells.Item(896, 7).Value = "Cilantro"
$ws.Cells.Item(896, 8).Value = "Sin especificar"
$ws.Cells.Item(896, 9).Value = "Primera"
$ws.Cells.Item(896, 10).Value = 70
$ws.Cells.Item(896, 11).Value = 8000
$ws.Cells.Item(896, 12).Value = 8000
$ws.Cells.Item(896, 13).Value = 8000
$ws.Cells.Item(896, 14).Value = "$/caja 36 atados"
$ws.Cells.Item(896, 15).Value = "Región Metropolitana"
$ws.Cells.Item(896, 16).Value = 222
$ws.Cells.Item(896, 17).Value = 36
$ws.Cells.Item(896, 18).Value = "Hortaliza"

# Populate the second new row (897) with the new record.
$ws.Cells.Item(897, 1).Value = 9
$ws.Cells.Item(897, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(897, 3).Value = "Metropolitana"
$ws.Cells.Item(897, 4).Value = 45041
$ws.Cells.Item(897, 5).Value = 13
$ws.Cells.Item(897, 6).Value = 100112040
$ws.Cells.Item(897, 7).Value = "Cilantro"
$ws.Cells.Item(897, 8).Value = "Sin especificar"
$ws.Cells.Item(897, 9).Value = "Primera"
$ws.Cells.Item(897, 10).Value = 160
$ws.Cells.Item(897, 11).Value = 10000
$ws.Cells.Item(897, 12).Value = 11000
$ws.Cells.Item(897, 13).Value = 10500
$ws.Cells.Item(897, 14).Value = "$/docena de atados"
$ws.Cells.Item(897, 15).Value = "Región Metropolitana"
$ws.Cells.Item(897, 16).Value = 3500
$ws.Cells.Item(897, 17).Value = 3
$ws.Cells.Item(897, 18).Value = "Hortaliza"
